$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "test"
$ws.Range("B2").Value = "test"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1
$ws.Range("H2").Value = "0, 1"
$ws.Range("K2").Value = "0, 1"
$ws.Range("N2").Value = "0, 1"
$ws.Range("R2").Value = 0
$ws.Range("W2").Value = 0.1

$ws.Range("A2").Select()
